$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Shift the five yearly-period columns (D:H) one column to the left and
# --- bring in a new year of data (1401/12, published 1402-02-23) in column H,
# --- dropping the oldest year (1396/12) that used to live in column D.

# Row 8: period headers ("12 ماهه منتهی به ....")
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# Row 9: publish dates
$ws.Range("D9").Value = "1399-04-19 (10)"
$ws.Range("E9").Value = "1400-04-25 (10)"
$ws.Range("F9").Value = "1401-04-18 (8)"
$ws.Range("G9").Value = "1402-02-23 (7)"
# H9 is a bare ISO-looking date ("1402-02-23") which Excel's Value setter
# would auto-convert to a date serial; round-trip it through a formula +
# paste-special-values so it lands back in the cell as literal text (and
# keeps the cell's existing style untouched).
$ws.Range("ZZ1").Formula = '="1402-02-23"'
$ws.Range("ZZ1").Copy()
$ws.Range("H9").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$ws.Range("ZZ1").ClearContents()
$excel.CutCopyMode = $false

# Row 11: فروش (Sales)
$ws.Range("D11").Value = 3637846
$ws.Range("E11").Value = 4378854
$ws.Range("F11").Value = 3197773
$ws.Range("G11").Value = 5353001
$ws.Range("H11").Value = 6237391

# Row 12: بهای تمام شده کالای فروش رفته
$ws.Range("D12").Value = -3264883
$ws.Range("E12").Value = -4071547
$ws.Range("F12").Value = -2816097
$ws.Range("G12").Value = -4836917
$ws.Range("H12").Value = -5495390

# Row 13: سود (زیان) ناخالص
$ws.Range("D13").Value = 372963
$ws.Range("E13").Value = 307307
$ws.Range("F13").Value = 381675
$ws.Range("G13").Value = 516084
$ws.Range("H13").Value = 742002

# Row 14: هزینه های عمومی, اداری و تشکیلاتی
$ws.Range("D14").Value = -42594
$ws.Range("E14").Value = -28928
$ws.Range("F14").Value = -43787
$ws.Range("G14").Value = -38914
$ws.Range("H14").Value = -39711

# Row 15: هزینه کاهش ارزش دریافتنی‌ها (هزینه استثنایی) -- stays "-" for every year, unchanged

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = 37649
$ws.Range("E16").Value = -65198
$ws.Range("F16").Value = 125421
$ws.Range("G16").Value = 56059
$ws.Range("H16").Value = 17054

# Row 17: سود (زیان) عملیاتی
$ws.Range("D17").Value = 368018
$ws.Range("E17").Value = 213181
$ws.Range("F17").Value = 463310
$ws.Range("G17").Value = 533229
$ws.Range("H17").Value = 719344

# Row 18: هزینه های مالی
$ws.Range("D18").Value = -632
$ws.Range("E18").Value = -710
$ws.Range("F18").Value = "-"
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = -10894

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 34662
$ws.Range("E19").Value = 89731
$ws.Range("F19").Value = 79305
$ws.Range("G19").Value = 132667
$ws.Range("H19").Value = 205547

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 402048
$ws.Range("E20").Value = 302202
$ws.Range("F20").Value = 542614
$ws.Range("G20").Value = 665896
$ws.Range("H20").Value = 913997

# Row 21: مالیات
$ws.Range("D21").Value = -52380
$ws.Range("E21").Value = -46864
$ws.Range("F21").Value = -52599
$ws.Range("G21").Value = -82618
$ws.Range("H21").Value = -63785

# Row 22: سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 349668
$ws.Range("E22").Value = 255338
$ws.Range("F22").Value = 490016
$ws.Range("G22").Value = 583278
$ws.Range("H22").Value = 850212

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی -- stays "-" for every year, unchanged

# Row 24: سود (زیان) خالص
$ws.Range("D24").Value = 349668
$ws.Range("E24").Value = 255338
$ws.Range("F24").Value = 490016
$ws.Range("G24").Value = 583278
$ws.Range("H24").Value = 850212

# Row 25: سود هر سهم پس از کسر مالیات -- all zero, unchanged

# Row 26: سرمایه
$ws.Range("D26").Value = 434967
$ws.Range("E26").Value = 584649
$ws.Range("F26").Value = 751871
$ws.Range("G26").Value = 1042246
$ws.Range("H26").Value = 779279

# Row 27: سود هر سهم بر اساس آخرین سرمایه -- all zero, unchanged
